# Trade #40 closed at 2026-02-17 15:28:49 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.24   # Current Capital
$summary.Range("B4").Value = 0.24      # Total P&L $
$summary.Range("B5").Value = 0.12      # Total P&L %
$summary.Range("B6").Value = 40        # Total Trades
$summary.Range("B7").Value = 12        # Winning Trades
$summary.Range("B9").Value = 30        # Win Rate %

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.24     # Capital (MarketMaking)
$status.Range("D4").Value = 40         # Trades
$status.Range("E4").Value = 0.24       # P&L $
$status.Range("F4").Value = 0.24       # P&L %
$status.Range("G4").Value = 30         # Win Rate %

# --- All Trades & MarketMaking sheets: close out trade #40 (row 41) ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G41").Value = 0.98
    $ws.Range("H41").Value = "CLOSED"
    $ws.Range("I41").Value = 104.1667
    $ws.Range("J41").Value = 0.5
    $ws.Range("K41").Value = 100.24
    $ws.Range("P41").Value = "early_exit"
    $ws.Range("Q41").Value = 5.09
}
